# Edit corresponding to the commit:
#   1) Slide 6's table (graphic frame) switches to a different built-in
#      table style ({B0C25C35-D6EF-4448-A2E9-1FF28AB147B7}).
#   2) The presentation's theme colour scheme changes from the
#      "Integral" palette to the standard "Office" palette (the deck's
#      Design/SlideMaster theme, persisted in ppt/theme/theme2.xml).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 6 --------------------------------------------
$slide  = $p.Slides.Item(6)
$shape  = $slide.Shapes.Item(2)          # the Google Shape;127;p18 graphicFrame/table
$table  = $shape.Table

# Table styles are applied via ApplyStyle, not by assigning .Style directly.
$table.ApplyStyle("{B0C25C35-D6EF-4448-A2E9-1FF28AB147B7}")

# --- 2. Theme colour scheme -------------------------------------------------
$master = $p.SlideMaster
$scheme = $master.ColorScheme

# RGB() packs as r + g*256 + b*65536 (standard COM colour order).
$scheme.Colors(1).RGB  = 0          # dk1      000000
$scheme.Colors(2).RGB  = 16777215   # lt1      FFFFFF
$scheme.Colors(3).RGB  = 6968388    # dk2      44546A
$scheme.Colors(4).RGB  = 15132391   # lt2      E7E6E6
$scheme.Colors(5).RGB  = 13998939   # accent1  5B9BD5
$scheme.Colors(6).RGB  = 3243501    # accent2  ED7D31
$scheme.Colors(7).RGB  = 10855845   # accent3  A5A5A5
$scheme.Colors(8).RGB  = 49407      # accent4  FFC000
$scheme.Colors(9).RGB  = 12874308   # accent5  4472C4
$scheme.Colors(10).RGB = 4697456    # accent6  70AD47
$scheme.Colors(11).RGB = 12673797   # hlink    0563C1
$scheme.Colors(12).RGB = 7491477    # folHlink 954F72
